$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.6 = 50809.83 pesos`n✅ 50809.83 pesos = 12.51 = 971.27 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the tasas rates on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 79.34999999999999
$wsTasas.Range("O10").Value = 4031.76
$wsTasas.Range("N12").Value = 4060
$wsTasas.Range("O12").Value = 77.61
